# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1.742940831014585)
    3  = @(0.04763786555579896, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.369736951971621)
    4  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    5  = @(0.04763786555579896, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.369736951971621)
    6  = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 0.496779210170732, 6.740334628841572)
    7  = @(0.127881588408715, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.094976487407548)
    8  = @(0.6753301551942219, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 2.290389397800092)
    9  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    10 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    11 = @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1.742940831014585)
    12 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    13 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    14 = @(0.6753301551942219, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.642425054193055)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
